# 🔄 Actualización automática del tracker
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in resultado/profit for rows whose matches have concluded
$ws.Range("G27").Value = "Acierto"
$ws.Range("H27").Value = 0.73

$ws.Range("G41").Value = "Fallo"
$ws.Range("H41").Value = -1

$ws.Range("G51").Value = "Acierto"
$ws.Range("H51").Value = 0.67

$ws.Range("G53").Value = "Acierto"
$ws.Range("H53").Value = 2

$ws.Range("G54").Value = "Acierto"
$ws.Range("H54").Value = 1

$ws.Range("G56").Value = "Fallo"
$ws.Range("H56").Value = -1

$ws.Range("G57").Value = "Fallo"
$ws.Range("H57").Value = -1

$ws.Range("G60").Value = "Fallo"
$ws.Range("H60").Value = -1

$ws.Range("G65").Value = "Fallo"
$ws.Range("H65").Value = -1

$ws.Range("G66").Value = "Fallo"
$ws.Range("H66").Value = -1

$ws.Range("G67").Value = "Fallo"
$ws.Range("H67").Value = -1

$ws.Range("G68").Value = "Acierto"
$ws.Range("H68").Value = 2.75

# Append new event row at the bottom of the tracker
$ws.Range("A71").Value = 14601351

$dateCell = $ws.Range("B71")
$dateCell.NumberFormat = "@"
$dateCell.Value = "2025-09-10"
$dateCell.ClearFormats()

$ws.Range("C71").Value = "Luca Castelnuovo"
$ws.Range("D71").Value = "Marat Sharipov"
$ws.Range("E71").Value = "Gana Luca Castelnuovo"
$ws.Range("F71").Value = 3

# Result / profit for this match are still pending -> leave blank (like row 70 before it was settled)
$ws.Range("G71").Value = ""
$ws.Range("G71").Font.Bold = $false
$ws.Range("H71").Value = ""
$ws.Range("H71").Font.Bold = $false
